$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All rows in column C (C2:C252) get their fitness value updated to 7310
$ws.Range("C2:C252").Value = 7310
